# Update 上海-漫展信息.xlsx per "output generated at 456a3b4" refresh.
# Sheets (1-indexed): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------
# Sheet 1 (展览): row10 content replaced (event cancelled -> new event),
# plus G46 flips from a sellable price to "不可售", plus many F (sold
# count) bumps.
# ---------------------------------------------------------------

$ws1.Range("C10").Value = "上海·无穹-中国 航天沉浸艺术展"
$ws1.Range("D10").Value = "上海浦东新区樱花路869号3F 上海喜玛拉雅美术馆"
$ws1.Range("E10").Value = "2024.06.08 10:00-10.07 20:00"
$ws1.Range("F10").Value = 0
$ws1.Range("G10").Value = 78
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=86957"
$ws1.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202406/Bus3lAnI1717558639134.jpeg"

$ws1.Range("G46").Value = "不可售"

$ws1.Range("F2").Value = 1767
$ws1.Range("F3").Value = 10332
$ws1.Range("F8").Value = 1679
$ws1.Range("F9").Value = 417
$ws1.Range("F11").Value = 232
$ws1.Range("F13").Value = 512
$ws1.Range("F15").Value = 137
$ws1.Range("F17").Value = 17
$ws1.Range("F19").Value = 112
$ws1.Range("F20").Value = 380
$ws1.Range("F23").Value = 32
$ws1.Range("F24").Value = 110
$ws1.Range("F25").Value = 1182
$ws1.Range("F26").Value = 205
$ws1.Range("F27").Value = 599
$ws1.Range("F31").Value = 587
$ws1.Range("F32").Value = 244
$ws1.Range("F35").Value = 738
$ws1.Range("F36").Value = 109206
$ws1.Range("F37").Value = 784
$ws1.Range("F38").Value = 530
$ws1.Range("F39").Value = 1272
$ws1.Range("F40").Value = 825
$ws1.Range("F41").Value = 718
$ws1.Range("F43").Value = 352
$ws1.Range("F44").Value = 35
$ws1.Range("F45").Value = 720

# ---------------------------------------------------------------
# Sheet 2 (演出): F (sold count) bumps + one refreshed cover image URL.
# ---------------------------------------------------------------

$ws2.Range("F6").Value = 92
$ws2.Range("F11").Value = 67
$ws2.Range("F18").Value = 1125
$ws2.Range("F20").Value = 1928
$ws2.Range("F21").Value = 1928
$ws2.Range("F22").Value = 1111
$ws2.Range("F23").Value = 340
$ws2.Range("F24").Value = 691
$ws2.Range("F36").Value = 196
$ws2.Range("F38").Value = 31
$ws2.Range("F40").Value = 135
$ws2.Range("F42").Value = 10
$ws2.Range("F43").Value = 8
$ws2.Range("I43").Value = "//i2.hdslb.com/bfs/openplatform/202406/heIkgSe71717657445122.jpeg"
$ws2.Range("F44").Value = 79

# ---------------------------------------------------------------
# Sheet 3 (本地生活): F (sold count) bumps, then a brand-new row 13 is
# appended (new local-life listing), extending the used range to I13.
# ---------------------------------------------------------------

$ws3.Range("F2").Value = 112
$ws3.Range("F4").Value = 824
$ws3.Range("F5").Value = 206
$ws3.Range("F6").Value = 2558
$ws3.Range("F7").Value = 4207
$ws3.Range("F10").Value = 389
$ws3.Range("F11").Value = 266
$ws3.Range("F12").Value = 238

# Append row 13, matching the style already used in column A (A2:A12).
$ws3.Range("A12").Copy($ws3.Range("A13"))
$ws3.Range("A13").Value = 12

$ws3.Range("B13").NumberFormat = "@"
$ws3.Range("B13").Value = "2024-06-14"

$ws3.Range("C13").Value = "上海·「排球少年!!垃圾场决战」主题店"
$ws3.Range("D13").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"

$ws3.Range("E13").NumberFormat = "@"
$ws3.Range("E13").Value = "2024.06.14 00:00-07.07 23:59"

$ws3.Range("F13").Value = 32
$ws3.Range("G13").Value = 10
$ws3.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=86948"
$ws3.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202406/GxQRrJ2V1717655942245.png"

# ---------------------------------------------------------------
# Sheet 4 (全部类型): combined-view F (sold count) bumps mirroring the
# other three sheets (values occasionally diverge slightly from their
# source sheet, matching the upstream export).
# ---------------------------------------------------------------

$ws4.Range("F2").Value = 1767
$ws4.Range("F4").Value = 824
$ws4.Range("F7").Value = 10333
$ws4.Range("F8").Value = 206
$ws4.Range("F9").Value = 4207
$ws4.Range("F11").Value = 389
$ws4.Range("F13").Value = 266
$ws4.Range("F14").Value = 1679
$ws4.Range("F15").Value = 417
$ws4.Range("F16").Value = 232
$ws4.Range("F19").Value = 137
$ws4.Range("F22").Value = 112
$ws4.Range("F23").Value = 1125
$ws4.Range("F24").Value = 380
$ws4.Range("F26").Value = 32
$ws4.Range("F27").Value = 1928
$ws4.Range("F28").Value = 1111
$ws4.Range("F29").Value = 1182
$ws4.Range("F34").Value = 587
$ws4.Range("F38").Value = 738
$ws4.Range("F41").Value = 784
$ws4.Range("F42").Value = 530
$ws4.Range("F43").Value = 825
$ws4.Range("F44").Value = 718
$ws4.Range("F46").Value = 352
$ws4.Range("F47").Value = 135
$ws4.Range("F48").Value = 720
$ws4.Range("F50").Value = 79

Write-Output "edit complete"
